$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D..P -> E..Q)
$ws.Columns("D").EntireColumn.Insert()

# Fill in the new "Year" boundary column
$ws.Range("D10").Value = "Year"
$ws.Range("D11").Value = 2030
$ws.Range("D12").Value = 2030

# Match the final selection left by the author
$ws.Range("K15").Select()
